$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
  'räumen',
  'jagen',
  'äußern',
  'boxen',
  'streichen',
  'sprengen',
  'brauchen',
  'kosten',
  'rufen',
  'kehren',
  'stammen',
  'bergen',
  'achten',
  'fallen',
  'knarren',
  'sinken',
  'wachsen',
  'filmen',
  'backen',
  'bauen',
  'schwören',
  'flehen',
  'spielen',
  'grüßen',
  'sorgen',
  'suchen',
  'wundern',
  'seufzen',
  'biegen',
  'werfen',
  'schrecken',
  'arten',
  'schreiten',
  'wirken',
  'warnen',
  'trauen',
  'bellen',
  'stecken',
  'lügen',
  'zünden',
  'pfeifen',
  'heulen',
  'klettern',
  'treiben',
  'formen',
  'betteln',
  'gründen',
  'schmecken',
  'hauen',
  'quälen',
  'dringen',
  'heilen',
  'decken',
  'spüren',
  'reizen',
  'wehtun',
  'pflanzen',
  'rasen',
  'töten',
  'heben',
  'sichern',
  'führen',
  'sterben',
  'kichern',
  'münzen',
  'ehren',
  'liefern',
  'fangen',
  'schwingen',
  'altern',
  'fahren',
  'jubeln',
  'freuen',
  'werden',
  'geben',
  'loben',
  'bitten',
  'folgen',
  'saufen',
  'spinnen',
  'irren',
  'helfen',
  'weichen',
  'schwächen',
  'schlucken',
  'schenken',
  'feiern',
  'bluten',
  'flüchten',
  'kümmern',
  'lockern',
  'dienen',
  'tollen',
  'klingen',
  'graben',
  'wüten',
  'zielen',
  'lesen',
  'mögen',
  'siegen',
  'zögern',
  'runden',
  'tropfen',
  'planen',
  'schulden',
  'gelten',
  'mauern',
  'fischen',
  'sperren',
  'stehlen',
  'scheitern',
  'fließen',
  'platzen',
  'ärgern',
  'scheinen',
  'zeigen',
  'wenden',
  'drehen',
  'trennen',
  'malen',
  'machen',
  'hören',
  'ändern',
  'erben',
  'liegen',
  'fällen',
  'greifen',
  'enden'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}

